$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column T (shift cells right) so it inherits column S's
# formatting exactly (same style indices as the existing 2022 column),
# mirroring the "add a 2023 column" edit described by the diff.
$ws.Columns("T").Insert(-4161)

# Row 1 grew slightly taller to fit the extra column header text.
$ws.Rows(1).RowHeight = 39.75

# Header year for the new column.
$ws.Range("T4").Value = 2023

# Data rows for 2023.
$ws.Range("T5").Value = 22.606300992622124
$ws.Range("T6").Value = 13.621194578764559
$ws.Range("T7").Value = 38.913029379337182
$ws.Range("T8").Value = 19.215987701767872
$ws.Range("T9").Value = "-"
$ws.Range("T10").Value = 19.9288256227758
$ws.Range("T11").Value = 48.820179007323027
$ws.Range("T12").Value = 18.458698661744346
$ws.Range("T13").Value = 7.704160246533128
$ws.Range("T14").Value = 29.197080291970806

Write-Output "Added 2023 column (T) with 11 data points"
